$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Proximity sheet: append rows 36-40
# ---------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "15:14:17", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:14:28", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:14:46", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:14:57", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:14:59", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 36
for ($i = 0; $i -lt $proximityRows.Length; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]

    # Column A holds a date-like string ("2026-02-01"); force it to stay as
    # plain text instead of being auto-converted into a date serial value.
    $cellA = $proximity.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $proximity.Cells.Item($r, 2).Value = $row[1]
    $proximity.Cells.Item($r, 3).Value = $row[2]
    $proximity.Cells.Item($r, 4).Value = $row[3]
    $proximity.Cells.Item($r, 5).Value = $row[4]
    $proximity.Cells.Item($r, 6).Value = $row[5]
}

# ---------------------------------------------------------------
# Camera sheet: append rows 8-11
# ---------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "15:14:16", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:14:30", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:14:45", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:14:59", "15:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 8
for ($i = 0; $i -lt $cameraRows.Length; $i++) {
    $r = $startRow + $i
    $row = $cameraRows[$i]

    $cellA = $camera.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $camera.Cells.Item($r, 2).Value = $row[1]
    $camera.Cells.Item($r, 3).Value = $row[2]
    $camera.Cells.Item($r, 4).Value = $row[3]
    $camera.Cells.Item($r, 5).Value = $row[4]
    $camera.Cells.Item($r, 6).Value = $row[5]
}
